$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp header in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 20:52"

# Row 16 (Canada): update B, C, E
$ws.Cells.Item(16, 2).Value = 24292
$ws.Cells.Item(16, 3).Value = 974
$ws.Cells.Item(16, 5).Value = 16473

# Row 25 (India): update B, C, D, E, G, H
$ws.Cells.Item(25, 2).Value = 9205
$ws.Cells.Item(25, 3).Value = 759
$ws.Cells.Item(25, 4).Value = 1080
$ws.Cells.Item(25, 5).Value = 7794
$ws.Cells.Item(25, 7).Value = 43
$ws.Cells.Item(25, 8).Value = 331

# Row 42 (Emiratos Arabes Unidos): update B, C, D, E, G, H
$ws.Cells.Item(42, 2).Value = 4123
$ws.Cells.Item(42, 3).Value = 387
$ws.Cells.Item(42, 4).Value = 680
$ws.Cells.Item(42, 5).Value = 3421
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(42, 8).Value = 22

# Row 108 (Estado de Palestina): update B, C, E
$ws.Cells.Item(108, 2).Value = 290
$ws.Cells.Item(108, 3).Value = 22
$ws.Cells.Item(108, 5).Value = 230

# Ruanda moves up, now ranks ahead of El Salvador and Camboya (new data),
# El Salvador and Camboya shift down one row each, keeping their own data.
$ws.Cells.Item(128, 1).Value = "Ruanda"
$ws.Cells.Item(128, 2).Value = 126
$ws.Cells.Item(128, 3).Value = 6
$ws.Cells.Item(128, 4).Value = 25
$ws.Cells.Item(128, 5).Value = 101
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 0

$ws.Cells.Item(129, 1).Value = "El Salvador"
$ws.Cells.Item(129, 2).Value = 125
$ws.Cells.Item(129, 3).Value = 7
$ws.Cells.Item(129, 4).Value = 21
$ws.Cells.Item(129, 5).Value = 98
$ws.Cells.Item(129, 6).Value = 3
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 6

$ws.Cells.Item(130, 1).Value = "Camboya"
$ws.Cells.Item(130, 2).Value = 122
$ws.Cells.Item(130, 3).Value = 2
$ws.Cells.Item(130, 4).Value = 77
$ws.Cells.Item(130, 5).Value = 45
$ws.Cells.Item(130, 6).Value = 1
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 0

# Polinesia Francesa moves up, now ranks ahead of Uganda and Islas Caimanes
# (new data); Uganda and Islas Caimanes shift down one row each, keeping
# their own data.
$ws.Cells.Item(143, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(143, 2).Value = 53
$ws.Cells.Item(143, 3).Value = 2
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 53
$ws.Cells.Item(143, 6).Value = 1
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 0

$ws.Cells.Item(144, 1).Value = "Uganda"
$ws.Cells.Item(144, 2).Value = 53
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 4
$ws.Cells.Item(144, 5).Value = 49
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 0

$ws.Cells.Item(145, 1).Value = "Islas Caimanes"
$ws.Cells.Item(145, 2).Value = 53
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 6
$ws.Cells.Item(145, 5).Value = 46
$ws.Cells.Item(145, 6).Value = 3
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 1
